$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column stores values as plain text (it can contain
# thousands-grouped numbers like "42.243.65"). When a new price looks
# like an ordinary decimal (e.g. "312.60") Excel would otherwise parse it
# straight into a number and drop the trailing zero, so we momentarily
# force Text format for the assignment and then restore General -
# matching how the cell was formatted before the edit.

# Row 2
$ws.Cells.Item(2, 4).Value = '42.243.65'
$ws.Cells.Item(2, 5).Value = '  -0.33%  '

# Row 3
$ws.Cells.Item(3, 4).Value = '2.279.68'
$ws.Cells.Item(3, 5).Value = '  -0.61%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.23%  '

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '312.60'
$ws.Cells.Item(5, 4).NumberFormat = 'General'
$ws.Cells.Item(5, 5).Value = '  -1.53%  '

# Row 6
$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '101.80'
$ws.Cells.Item(6, 4).NumberFormat = 'General'
$ws.Cells.Item(6, 5).Value = '  -0.61%  '

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '0.614'
$ws.Cells.Item(7, 4).NumberFormat = 'General'
$ws.Cells.Item(7, 5).Value = '  -2.36%  '

# Row 8
$ws.Cells.Item(8, 5).Value = '  -0.12%  '

# Row 9
$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.595'
$ws.Cells.Item(9, 4).NumberFormat = 'General'
$ws.Cells.Item(9, 5).Value = '  -1.74%  '

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '38.72'
$ws.Cells.Item(10, 4).NumberFormat = 'General'
$ws.Cells.Item(10, 5).Value = '  -2.15%  '

# Row 11
$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.0897'
$ws.Cells.Item(11, 4).NumberFormat = 'General'
$ws.Cells.Item(11, 5).Value = '  -0.95%  '

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = '@'
$ws.Cells.Item(12, 4).Value = '8.18'
$ws.Cells.Item(12, 4).NumberFormat = 'General'
$ws.Cells.Item(12, 5).Value = '  -2.56%  '

# Row 13
$ws.Cells.Item(13, 5).Value = '  +1.26%  '

# Row 14
$ws.Cells.Item(14, 5).Value = '  +1.72%  '

# Row 15
$ws.Cells.Item(15, 5).Value = '  -1.63%  '

# Row 16
$ws.Cells.Item(16, 4).Value = '2.625.58'
$ws.Cells.Item(16, 5).Value = '  -0.62%  '

# Row 17
$ws.Cells.Item(17, 4).Value = '2.283.68'
$ws.Cells.Item(17, 5).Value = '  -0.29%  '

# Row 18
$ws.Cells.Item(18, 4).Value = '42.202.87'
$ws.Cells.Item(18, 5).Value = '  -0.27%  '

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = '@'
$ws.Cells.Item(19, 4).Value = '7.22'
$ws.Cells.Item(19, 4).NumberFormat = 'General'
$ws.Cells.Item(19, 5).Value = '  -2.30%  '

# Row 20
$ws.Cells.Item(20, 5).Value = '  -1.05%  '

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '13.43'
$ws.Cells.Item(21, 4).NumberFormat = 'General'
$ws.Cells.Item(21, 5).Value = '  +8.14%  '

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '72.80'
$ws.Cells.Item(22, 4).NumberFormat = 'General'
$ws.Cells.Item(22, 5).Value = '  -0.82%  '

# Row 23
$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '3.50'
$ws.Cells.Item(23, 4).NumberFormat = 'General'
$ws.Cells.Item(23, 5).Value = '  -1.19%  '

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '263.10'
$ws.Cells.Item(24, 4).NumberFormat = 'General'
$ws.Cells.Item(24, 5).Value = '  -4.57%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  -4.19%  '

# Row 26
$ws.Cells.Item(26, 5).Value = '  +0.30%  '

# Row 27
$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '10.57'
$ws.Cells.Item(27, 4).NumberFormat = 'General'
$ws.Cells.Item(27, 5).Value = '  -2.22%  '

# Row 28
$ws.Cells.Item(28, 2).Value = 'Filecoin'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '6.88'
$ws.Cells.Item(28, 4).NumberFormat = 'General'
$ws.Cells.Item(28, 5).Value = '  +13.56%  '

# Row 29
$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '2.27'
$ws.Cells.Item(29, 4).NumberFormat = 'General'
$ws.Cells.Item(29, 5).Value = '  -3.99%  '

# Row 30
$ws.Cells.Item(30, 5).Value = '  -1.69%  '

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '35.86'
$ws.Cells.Item(31, 4).NumberFormat = 'General'
$ws.Cells.Item(31, 5).Value = '  -4.35%  '

# Row 32
$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '164.04'
$ws.Cells.Item(32, 4).NumberFormat = 'General'
$ws.Cells.Item(32, 5).Value = '  -1.10%  '

# Row 33
$ws.Cells.Item(33, 5).Value = '  -1.80%  '

# Row 34
$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '0.130'
$ws.Cells.Item(34, 4).NumberFormat = 'General'
$ws.Cells.Item(34, 5).Value = '  -2.81%  '

# Row 35
$ws.Cells.Item(35, 5).Value = '  -2.14%  '

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '4.47'
$ws.Cells.Item(37, 4).NumberFormat = 'General'
$ws.Cells.Item(37, 5).Value = '  -1.97%  '

# Row 38
$ws.Cells.Item(38, 5).Value = '  -4.27%  '

# Row 39
$ws.Cells.Item(39, 5).Value = '  -0.33%  '

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = '@'
$ws.Cells.Item(40, 4).Value = '2.64'
$ws.Cells.Item(40, 4).NumberFormat = 'General'
$ws.Cells.Item(40, 5).Value = '  -3.63%  '

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '1.57'
$ws.Cells.Item(41, 4).NumberFormat = 'General'
$ws.Cells.Item(41, 5).Value = '  +5.28%  '

# Row 42
$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '97.48'
$ws.Cells.Item(42, 4).NumberFormat = 'General'
$ws.Cells.Item(42, 5).Value = '  +1.80%  '

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = '@'
$ws.Cells.Item(43, 4).Value = '68.92'
$ws.Cells.Item(43, 4).NumberFormat = 'General'
$ws.Cells.Item(43, 5).Value = '  -0.98%  '

# Row 44
$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '0.226'
$ws.Cells.Item(44, 4).NumberFormat = 'General'
$ws.Cells.Item(44, 5).Value = '  +0.66%  '

# Row 45
$ws.Cells.Item(45, 5).Value = '  -0.50%  '

# Row 46
$ws.Cells.Item(46, 5).Value = '  -0.63%  '

# Row 47
$ws.Cells.Item(47, 4).Value = '1.698.68'
$ws.Cells.Item(47, 5).Value = '  +6.24%  '

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '78.71'
$ws.Cells.Item(48, 4).NumberFormat = 'General'
$ws.Cells.Item(48, 5).Value = '  -0.50%  '

# Row 49
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '109.71'
$ws.Cells.Item(49, 4).NumberFormat = 'General'
$ws.Cells.Item(49, 5).Value = '  -2.17%  '

# Row 50
$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '5.19'
$ws.Cells.Item(50, 4).NumberFormat = 'General'
$ws.Cells.Item(50, 5).Value = '  -0.97%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -3.72%  '
